$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A45").Value = "$ 27.386 CLP 29-10-20"
$ws.Range("A46").Value = "$ 27.386 CLP 29-10-20"
